# Update "plan de test" slide: reword / split a few bullet lines in the
# "Forme libre : forme 4" shape (Gestion de profil test plan) and grow the
# shape a bit to make room for the extra lines.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)
$shape = $grp.GroupItems.Item(2)

# ---------------------------------------------------------------------
# 1) Resize / reposition the shape (child-space off/ext changes from
#    x=178838,y=2392897,cx=5810680,cy=2072204 to
#    x=204905,y=2185511,cx=5810680,cy=2351854). This host's Shape
#    Left/Top/Height setters write straight into the shape's own <a:xfrm>
#    (point-value * 12700 EMU/pt) without re-applying the parent group's
#    off/ext transform, so we feed them the raw target EMU values
#    converted back to points (with a tiny epsilon on a couple of them
#    to dodge float-truncation landing one EMU short).
# ---------------------------------------------------------------------
$shape.Left = 16.134251968503936
$shape.Top = 172.08749031496063
$shape.Height = 185.18536433070867

# ---------------------------------------------------------------------
# 2) Text edits on the paragraphs of the shape's TextFrame.
#    Paragraphs (1-based) before editing:
#      1 Scope
#      2 Entrer dans un champ type email
#      3 Entrer dans un champ type charactères
#      4 Entrer dans champ type numérique
#      5 Suppression du profil
#      6 (empty)
#      7 (tab)
# ---------------------------------------------------------------------
$tr = $shape.TextFrame.TextRange

# --- paragraph 2: "...type email" -> "...type email pour l'email" ---
$para = $tr.Paragraphs(2, 1)
$para.Text = "__tmp__"
$para = $tr.Paragraphs(2, 1)
$para.Text = "Entrer dans un champ type email pour l’email"

# New paragraph right after it.
$para = $tr.Paragraphs(2, 1)
$para.InsertAfter("`rEntrer dans un champ type non-email pour l’email")

# --- paragraph 4 (was 3): "...type charactères " -> "...pour l'adresse et lieu" ---
$para = $tr.Paragraphs(4, 1)
$para.Text = "__tmp__"
$para = $tr.Paragraphs(4, 1)
$para.Text = "Entrer dans un champ type charactères pour l’adresse et lieu"

# New paragraph right after it.
$para = $tr.Paragraphs(4, 1)
$para.InsertAfter("`rEntrer dans un champ type non-charactères pour l’adresse et lieu")

# Another new paragraph right after that one.
$para = $tr.Paragraphs(5, 1)
$para.InsertAfter("`rEntrer dans champ type numérique pour le NPA")

# --- paragraph 7 (was 4, original "...type numérique"): -> "...non-numérique pour le NPA" ---
$para = $tr.Paragraphs(7, 1)
$para.Text = "__tmp__"
$para = $tr.Paragraphs(7, 1)
$para.Text = "Entrer dans champ type non-numérique pour le NPA"
